$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Ciudades")

# Update the "last updated" timestamp string in A1 (08:20 -> 08:50)
$ws.Range("A1").Value = "Datos actualizados a 31 de Marzo de 2020 a las 08:50"

# Update Recuperados (D) and Muertes (E) values for the affected rows
# Row 17 - A Coruña
$ws.Range("D17").Value = 1250
$ws.Range("E17").Value = 51

# Row 29 - Pontevedra
$ws.Range("D29").Value = 998
$ws.Range("E29").Value = 16

# Row 46 - Ourense
$ws.Range("D46").Value = 411
$ws.Range("E46").Value = 12

# Row 49 - Lugo
$ws.Range("D49").Value = 241
$ws.Range("E49").Value = 7
